$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new job link in A4, matching style of existing hyperlink cells (A2, A3)
$ws.Range("A4").Value = "https://www.bestjobs.eu/ro/"
$ws.Hyperlinks.Add($ws.Range("A4"), "https://www.bestjobs.eu/ro/", "", "", "https://www.bestjobs.eu/ro/") | Out-Null
$ws.Range("A4").Style = "Hyperlink"

# Update the selection to A5 like in the diff
$ws.Range("A5").Select()
